$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed-date) column C for rows 2-7 from 45204 to 45207
$ws.Range("C2:C7").Value = 45207
